$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.239.06"
$ws.Range("E2").Value = "  -0.78%  "
$ws.Range("D3").Value = "1.663.98"
$ws.Range("E3").Value = "  -0.37%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5235"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2647"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.77%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06286"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.85"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07765"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.46%  "
$ws.Range("D12").Value = "1.683.61"
$ws.Range("E12").Value = "  +0.73%  "
$ws.Range("E13").Value = "  +0.38%  "
$ws.Range("D14").Value = "1.890.70"
$ws.Range("E14").Value = "  -0.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5464"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.15%  "
$ws.Range("D16").Value = "0.0₅8162"
$ws.Range("E16").Value = "  -1.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.00"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.70%  "
$ws.Range("D18").Value = "26.258.42"
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.616"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "192.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.03"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.90%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.028"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.004"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "138.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1237"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.270"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.67%  "
$ws.Range("E28").Value = "  -0.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.417"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05976"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.276"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.544"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.273"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.582"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9606"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.26%  "
$ws.Range("E36").Value = "  -0.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5685"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.30%  "
$ws.Range("E39").Value = "  -0.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.972"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8527"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.001"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.52"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.16%  "
$ws.Range("D44").Value = "1.005.09"
$ws.Range("E44").Value = "  -7.88%  "
$ws.Range("D45").Value = "1.805.44"
$ws.Range("E45").Value = "  -0.41%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.74"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.12%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₈107"
$ws.Range("E47").Value = "  -1.62%  "
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.055"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4346"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05151"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.03%  "

Write-Host "Applied cryptos update"
